# fix: update and clean ups
# Adds a new "Narration" column (I) to the batch sheet with sample values,
# matching the style of the neighbouring AccountType column, and moves the
# active selection to K7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Narration" header + per-row values in column I
$ws.Range("I1").Value = "Narration"
$ws.Range("I2").Value = "test"
$ws.Range("I3").Value = "Test again"
$ws.Range("I4").Value = "Awesome test"

# Match formatting of the existing text-formatted columns (e.g. H)
$ws.Range("I1:I4").NumberFormat = "@"

# Update the sheet's active selection
$ws.Range("K7").Select() | Out-Null
